# Daily attendance processing
# Normalizes the "Recorded By" column (G) so that the automated "System"
# recorder entry is always listed last among the recorders for a session,
# while the remaining (human / other) recorder names keep a stable,
# deterministic (descending) relative order.

function Sort-Descending($arr) {
    $n = $arr.Length
    for ($i = 0; $i -lt $n; $i++) {
        for ($j = 0; $j -lt $n - $i - 1; $j++) {
            if ($arr[$j].CompareTo($arr[$j + 1]) -lt 0) {
                $tmp = $arr[$j]
                $arr[$j] = $arr[$j + 1]
                $arr[$j + 1] = $tmp
            }
        }
    }
    return $arr
}

function Reorder-Recorders($value) {
    if ($null -eq $value) {
        return $value
    }

    $parts = $value -split ","
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    # Count exact ("System", case-sensitive) occurrences and keep everything
    # else (including lowercase "system") in $rest.
    $systemCount = 0
    $rest = @()
    foreach ($p in $trimmed) {
        if ($p.Equals("System")) {
            $systemCount++
        } else {
            $rest += $p
        }
    }

    if ($systemCount -eq 0) {
        return $value
    }

    $rest = Sort-Descending $rest

    $result = @()
    foreach ($r in $rest) {
        $result += $r
    }
    for ($k = 0; $k -lt $systemCount; $k++) {
        $result += "System"
    }

    return ($result -join ", ")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$recordedByCol = 7  # Column G: "Recorded By"

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $recordedByCol)
    $current = $cell.Value()

    if ($current -is [string] -and $current.Length -gt 0) {
        $updated = Reorder-Recorders $current
        if (-not $updated.Equals($current)) {
            $cell.Value = $updated
        }
    }
}
